$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.586902141571045
$ws.Range("B1").Value = 1.682878851890564
$ws.Range("C1").Value = 1.733575344085693
$ws.Range("D1").Value = 2.305009603500366
$ws.Range("E1").Value = 3.355896949768066
